# "Varios fixes realizados y mensajes de usuario"
#
# The "Trimestre"/"Semestre" header pair in the Atenciones sheet is replaced
# by a single "Periodo" header: column L ("Semestre") is removed entirely
# (shifting every column to its right one position to the left), and the
# surviving column K header text ("Trimestre") is renamed to "Periodo".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "Semestre" column (L) - cells to its right shift left.
$ws.Range("L1").EntireColumn.Delete() | Out-Null

# The old "Trimestre" header (now in column K) becomes "Periodo".
$ws.Range("K5").Value = "Periodo"

# Restore the user's last selection as seen in the saved workbook.
$ws.Range("L19").Select() | Out-Null
